$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.503.26'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.063.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.43%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.056.61'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.39%  '
$ws.Range("E9").Value = '  +3.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.25'
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = '  +2.18%  '
$ws.Range("E12").Value = '  +4.72%  '
$ws.Range("E13").Value = '  +5.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.566.41'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.508.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.065.94'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.64%  '
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("E19").Value = '  +4.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '478.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.29%  '
$ws.Range("E22").Value = '  +2.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.24%  '
$ws.Range("E25").Value = '  +7.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("E27").Value = '  +5.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.24%  '
$ws.Range("E29").Value = '  +9.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.42%  '
$ws.Range("E33").Value = '  +8.87%  '
$ws.Range("E34").Value = '  +6.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.92'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.99'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '467.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0816'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.138.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0395'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.56%  '
$ws.Range("E41").Value = '  +4.24%  '
$ws.Range("E42").Value = '  +4.20%  '
$ws.Range("E43").Value = '  +8.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.95'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.05%  '
$ws.Range("E45").Value = '  +4.81%  '
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.47%  '
$ws.Range("E48").Value = '  +2.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0513'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '116.20'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.36%  '
